$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of column J into the matching rows of the new column K
# (only the rows that actually carry an explicit cell in column J) so the
# appended column inherits the same borders/fonts/number formats.
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("J7").Copy() | Out-Null
$ws.Range("K7").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").Copy() | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null

# Populate the new column K values (2020 data).
$ws.Range("K4").Value = 2020
$ws.Range("K6").Value = 5.9
$ws.Range("K7").Value = 1.5
$ws.Range("K8").Value = "-"

# Restore the selection seen in the saved workbook.
$ws.Range("L16").Select() | Out-Null
